# Adds measure M35 "Het project hanteert een agile architectuuraanpak" to the
# ICTU Kwaliteitsaanpak deck. A new slide is inserted right before the
# existing "M10: Het project kent een wekelijks projectoverleg" slide (slide
# 20), pushing that slide (and everything after it) down by one position.
#
# Implementation: duplicate the M10 slide (which uses the same simple
# Title + free-floating TextBox layout we need for the new slide), which
# places the duplicate immediately after the original. Then overwrite the
# *original* slide (still at position 20) with the new M35 title/body text,
# leaving the duplicate (now at position 21) holding the untouched M10
# content. Every subsequent slide keeps its content and simply shifts down
# by one position, which is exactly the desired outcome.

$p = $ppt.ActivePresentation

$m10Index = 20
$m10Slide = $p.Slides.Item($m10Index)

# Duplicate keeps the original's content/formatting intact and inserts the
# copy right after it (position 21).
$m10Slide.Duplicate() | Out-Null

# The shape at position 20 is still the original; repurpose it for M35.
$m35Slide = $p.Slides.Item($m10Index)
$m35Slide.Shapes.Item(1).TextFrame.TextRange.Text = "M35: Het project hanteert een agile architectuuraanpak"
$m35Slide.Shapes.Item(2).TextFrame.TextRange.Text = "Tijdens de voorfase verwerkt het project de door de opdrachtgever opgestelde projectstartarchitectuur (PSA) in een eerste versie van het softwarearchitectuurdocument (SAD). Tijdens de realisatiefase werkt het project het SAD bij op basis van nieuwe inzichten."
